$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-VQUG-001"
$ws.Range("F2").Value = "RMA-VQUG-1-1"
$ws.Range("J2").Value = "a7s5f000000xLcDAAU"

$ws.Range("E3").Value = "RMA-VQUG-002"
$ws.Range("F3").Value = "RMA-VQUG-1-2"
$ws.Range("J3").Value = "a7s5f000000xLcEAAU"

$ws.Range("E4").Value = "RMA-VQUG-003"
$ws.Range("F4").Value = "RMA-VQUG-1-3"
$ws.Range("J4").Value = "a7s5f000000xLcFAAU"
